$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: rows 1533-1538, column R: inlineStr (blank) -> numeric 0 ----
for ($r = 1533; $r -le 1538; $r++) {
    $ws.Cells.Item($r, 18).Value = 0
}

# ---- Step 2: append new historical rows 1539-1555 (columns A-Q); column R stays blank ----
$csv = @'
45596,4400.89990234375,4400.89990234375,4258.5,4299.4501953125,4299.4501953125,520292,2024,10,31,0,0,0,44,0,2,0
45597,4340,4359.14990234375,4310,4322.85009765625,4322.85009765625,55312,2024,11,1,0,0,0,44,0,0,0
45600,4342.0498046875,4342.0498046875,4227,4283.39990234375,4283.39990234375,560192,2024,11,4,0,0,0,45,0,0,0
45601,4294,4328,4182.5,4294,4294,859564,2024,11,5,0,0,0,45,0,0,0
45602,4300,4509,4300,4498.4501953125,4498.4501953125,618525,2024,11,6,0,0,0,45,0,0,0
45603,4509.60009765625,4539,4458,4517.14990234375,4517.14990234375,276984,2024,11,7,0,0,0,45,0,0,0
45604,4529,4552.89990234375,4460.0498046875,4484.5498046875,4484.5498046875,265978,2024,11,8,0,0,0,45,1,0,0
45607,4484.5498046875,4500,4392.64990234375,4450.0498046875,4450.0498046875,297456,2024,11,11,0,0,0,46,0,0,0
45608,4474,4540,4380.0498046875,4397.5,4397.5,477903,2024,11,12,0,0,0,46,0,0,0
45609,4360.9501953125,4372.10009765625,4220,4249.14990234375,4249.14990234375,364377,2024,11,13,0,0,0,46,0,0,0
45610,4249.0498046875,4360,4248.0498046875,4278.75,4278.75,393645,2024,11,14,0,0,0,46,0,0,0
45614,4276,4278.75,4167.0498046875,4182.35009765625,4182.35009765625,684781,2024,11,18,0,0,0,47,0,0,0
45615,4182,4229,4151.75,4190.14990234375,4190.14990234375,455486,2024,11,19,0,0,0,47,0,0,0
45617,4200,4237.60009765625,4126.60009765625,4214.10009765625,4214.10009765625,873511,2024,11,21,0,0,0,47,0,0,0
45618,4224.60009765625,4240,4185.10009765625,4223.25,4223.25,472525,2024,11,22,0,0,0,47,0,0,0
45621,4304.9501953125,4373.35009765625,4272.2001953125,4343.4501953125,4343.4501953125,704527,2024,11,25,0,0,0,48,0,0,0
45622,4375.5,4375.5,4232,4238.4501953125,4238.4501953125,269426,2024,11,26,0,0,0,48,0,0,0
'@

$lines = $csv -split "`n"
$data = New-Object 'object[,]' 17,17
for ($i = 0; $i -lt $lines.Length; $i++) {
    $fields = $lines[$i].Split(",")
    for ($j = 0; $j -lt $fields.Length; $j++) {
        $data[$i, $j] = [double]$fields[$j]
    }
}

$ws.Range("A1539:Q1555").Value = $data

# Column A on the new rows carries the same datetime display format as the rest of column A
$ws.Range("A1539:A1555").NumberFormat = "YYYY-MM-DD HH:MM:SS"

